# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The worksheet lists arrears ("Estado de Cuenta") rows for workers. The
# previous data (row 16 = KAREN SOFIA LEZAMA GARCES / CC 32937423, rows
# 17-42 = GLADYS DANIELA MARIN RIOS / CC 1047429612 for periods 2111..2401
# ascending) is replaced with an updated data set: GLADYS's periods are
# now listed descending (2401..2111) in rows 16-41 with refreshed
# "Valor Mora" / "Salario Basico" amounts, and KAREN's single row moves
# to the bottom (row 42), keeping her original values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 16-41: CC, 1047429612, GLADYS DANIELA MARIN RIOS, period, ValorMora(F), SalarioBasico(G)
$gladysRows = @(
    @{ Row = 16; Periodo = "2401"; Valor = 40000;  Salario = 1000000 },
    @{ Row = 17; Periodo = "2312"; Valor = 40000;  Salario = 1000000 },
    @{ Row = 18; Periodo = "2311"; Valor = 40000;  Salario = 1000000 },
    @{ Row = 19; Periodo = "2310"; Valor = 40000;  Salario = 1000000 },
    @{ Row = 20; Periodo = "2309"; Valor = 40000;  Salario = 1000000 },
    @{ Row = 21; Periodo = "2308"; Valor = 40000;  Salario = 1000000 },
    @{ Row = 22; Periodo = "2307"; Valor = 40000;  Salario = 1000000 },
    @{ Row = 23; Periodo = "2306"; Valor = 40000;  Salario = 1000000 },
    @{ Row = 24; Periodo = "2305"; Valor = 40000;  Salario = 1000000 },
    @{ Row = 25; Periodo = "2304"; Valor = 40000;  Salario = 1000000 },
    @{ Row = 26; Periodo = "2303"; Valor = 40000;  Salario = 1000000 },
    @{ Row = 27; Periodo = "2302"; Valor = 40000;  Salario = 1000000 },
    @{ Row = 28; Periodo = "2301"; Valor = 40000;  Salario = 1000000 },
    @{ Row = 29; Periodo = "2212"; Valor = 40000;  Salario = 1000000 },
    @{ Row = 30; Periodo = "2211"; Valor = 40000;  Salario = 1000000 },
    @{ Row = 31; Periodo = "2210"; Valor = 40000;  Salario = 1000000 },
    @{ Row = 32; Periodo = "2209"; Valor = 40000;  Salario = 1000000 },
    @{ Row = 33; Periodo = "2208"; Valor = 40000;  Salario = 1000000 },
    @{ Row = 34; Periodo = "2206"; Valor = 36341;  Salario = 1000000 },
    @{ Row = 35; Periodo = "2205"; Valor = 36341;  Salario = 1000000 },
    @{ Row = 36; Periodo = "2204"; Valor = 36341;  Salario = 1000000 },
    @{ Row = 37; Periodo = "2203"; Valor = 36341;  Salario = 1000000 },
    @{ Row = 38; Periodo = "2202"; Valor = 36341;  Salario = 1000000 },
    @{ Row = 39; Periodo = "2201"; Valor = 36341;  Salario = 1000000 },
    @{ Row = 40; Periodo = "2112"; Valor = 36341;  Salario = 1000000 },
    @{ Row = 41; Periodo = "2111"; Valor = 32707;  Salario = 1000000 }
)

foreach ($item in $gladysRows) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = "1047429612"
    $ws.Cells.Item($r, 4).Value = "GLADYS DANIELA MARIN RIOS"
    $ws.Cells.Item($r, 5).Value = $item.Periodo
    $ws.Cells.Item($r, 6).Value = $item.Valor
    $ws.Cells.Item($r, 7).Value = $item.Salario
}

# Row 42 (last data row / totals-adjacent row) now holds KAREN's original record.
$ws.Cells.Item(42, 3).Value = "32937423"
$ws.Cells.Item(42, 4).Value = "KAREN SOFIA LEZAMA GARCES"
$ws.Cells.Item(42, 5).Value = "1810"
$ws.Cells.Item(42, 6).Value = 32000
$ws.Cells.Item(42, 7).Value = 828116
